$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.789.63'
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').Value = '1.891.12'
$ws.Range('E3').Value = '  +1.18%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9997'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.59'
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4740'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2929'
$ws.Range('E8').Value = '  +0.41%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06540'
$ws.Range('E9').Value = '  +0.69%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.16'
$ws.Range('E10').Value = '  +0.39%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07803'
$ws.Range('E11').Value = '  +1.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '97.01'
$ws.Range('E12').Value = '  -0.59%  '
$ws.Range('D13').Value = '1.889.71'
$ws.Range('E13').Value = '  +0.93%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7391'
$ws.Range('E14').Value = '  -0.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.253'
$ws.Range('E15').Value = '  +2.39%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '285.61'
$ws.Range('E16').Value = '  +4.36%  '
$ws.Range('D17').Value = '30.767.46'
$ws.Range('E17').Value = '  +1.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.28'
$ws.Range('E18').Value = '  -0.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007558'
$ws.Range('E19').Value = '  +0.28%  '
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').Value = '2.136.78'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.334'
$ws.Range('E22').Value = '  +2.14%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.264'
$ws.Range('E24').Value = '  +1.50%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.234'
$ws.Range('E25').Value = '  -0.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.98'
$ws.Range('E26').Value = '  +0.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.02'
$ws.Range('E27').Value = '  +1.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.920'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.343'
$ws.Range('E29').Value = '  -1.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09754'
$ws.Range('E30').Value = '  -2.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.496'
$ws.Range('E31').Value = '  -0.23%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.304'
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.187'
$ws.Range('E33').Value = '  +1.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04867'
$ws.Range('E34').Value = '  +0.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.128'
$ws.Range('E35').Value = '  +0.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6977'
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.724'
$ws.Range('E37').Value = '  +0.35%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01893'
$ws.Range('E38').Value = '  +1.85%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.805'
$ws.Range('E39').Value = '  +2.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '76.26'
$ws.Range('E40').Value = '  +4.19%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.331'
$ws.Range('E41').Value = '  +0.36%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.996'
$ws.Range('E42').Value = '  +1.51%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4291'
$ws.Range('E43').Value = '  +2.22%  '
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8359'
$ws.Range('E45').Value = '  +0.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '101.79'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.537'
$ws.Range('E47').Value = '  +3.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.061'
$ws.Range('E48').Value = '  +0.88%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '35.64'
$ws.Range('E49').Value = '  +0.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '915.70'
$ws.Range('E50').Value = '  -1.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05755'
